# Update Sheets via scheduled runner: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) for the affected leve rows across the crafting-job worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H86").Value = 1931.6666
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 1997.5
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 1997.5
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -4243.5
$ws.Range("H89").Value = 1931.6666
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 1997.5
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 9987.5
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -21219.5
$ws.Range("H106").Value = 2246.2144
$ws.Range("I106").Value = 2412.75
$ws.Range("K106").Value = 2412.75
$ws.Range("M106").Value = -1781.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 927442.3
$ws.Range("I2").Value = 1390163.8
$ws.Range("K2").Value = 1390163.8
$ws.Range("M2").Value = -1390050.8
$ws.Range("H5").Value = 248.5
$ws.Range("I5").Value = 248.5
$ws.Range("K5").Value = 248.5
$ws.Range("M5").Value = -136.5
$ws.Range("H32").Value = 4144
$ws.Range("I32").Value = 3359.98
$ws.Range("K32").Value = 3359.98
$ws.Range("M32").Value = -3072.98
$ws.Range("H116").Value = 927442.3
$ws.Range("I116").Value = 1390163.8
$ws.Range("K116").Value = 1390163.8
$ws.Range("M116").Value = -1387869.8
$ws.Range("H132").Value = 1636.8918
$ws.Range("I132").Value = 1336.1724
$ws.Range("K132").Value = 4008.5172
$ws.Range("M132").Value = -1478.5172

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 927442.3
$ws.Range("I3").Value = 1390163.8
$ws.Range("K3").Value = 1390163.8
$ws.Range("M3").Value = -1390049.8
$ws.Range("H4").Value = 248.5
$ws.Range("I4").Value = 248.5
$ws.Range("K4").Value = 248.5
$ws.Range("M4").Value = -133.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H134").Value = 4405.8857
$ws.Range("I134").Value = 4809.643
$ws.Range("K134").Value = 14428.929
$ws.Range("M134").Value = -11893.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3208.125
$ws.Range("I86").Value = 3208.125
$ws.Range("K86").Value = 3208.125
$ws.Range("M86").Value = -2085.125
$ws.Range("H89").Value = 3208.125
$ws.Range("I89").Value = 3208.125
$ws.Range("K89").Value = 16040.625
$ws.Range("M89").Value = -10424.625
$ws.Range("H99").Value = 1112885.5
$ws.Range("I99").Value = 2501742.2
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 2501742.2
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -2500244.2
$ws.Range("N99").Value = -4796
$ws.Range("H122").Value = 3803.8333
$ws.Range("I122").Value = 2737
$ws.Range("K122").Value = 8211
$ws.Range("M122").Value = -5761
$ws.Range("H126").Value = 1112885.5
$ws.Range("I126").Value = 2501742.2
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 7505226.600000001
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -7502756.600000001
$ws.Range("N126").Value = -10340
$ws.Range("H132").Value = 2949.4375
$ws.Range("I132").Value = 1818.3
$ws.Range("K132").Value = 5454.9
$ws.Range("M132").Value = -2924.9
$ws.Range("H134").Value = 1959.0625
$ws.Range("I134").Value = 1965.5714
$ws.Range("K134").Value = 5896.7142
$ws.Range("M134").Value = -3361.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 10778.444
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 10778.444
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 32335.332
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -32673.332
$ws.Range("H40").Value = 108.333336
$ws.Range("I40").Value = 97.5
$ws.Range("K40").Value = 390
$ws.Range("M40").Value = -321
$ws.Range("H75").Value = 500
$ws.Range("J75").Value = 500
$ws.Range("L75").Value = 1500
$ws.Range("N75").Value = -3496
$ws.Range("H78").Value = 500
$ws.Range("J78").Value = 500
$ws.Range("L78").Value = 4500
$ws.Range("N78").Value = -14484
$ws.Range("H131").Value = 11872.069
$ws.Range("J131").Value = 12710.881
$ws.Range("L131").Value = 38132.643
$ws.Range("N131").Value = -48212.643

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4098705.8
$ws.Range("J7").Value = 1035599.6
$ws.Range("L7").Value = 1035599.6
$ws.Range("N7").Value = -1035823.6
$ws.Range("H8").Value = 4098705.8
$ws.Range("J8").Value = 1035599.6
$ws.Range("L8").Value = 1035599.6
$ws.Range("N8").Value = -1035877.6
$ws.Range("H11").Value = 5355497
$ws.Range("I11").Value = 6883300
$ws.Range("J11").Value = 2855456.2
$ws.Range("K11").Value = 6883300
$ws.Range("L11").Value = 2855456.2
$ws.Range("M11").Value = -6883161
$ws.Range("N11").Value = -2855734.2
$ws.Range("H14").Value = 2025376.9
$ws.Range("I14").Value = 3800000
$ws.Range("K14").Value = 3800000
$ws.Range("M14").Value = -3799832
$ws.Range("H18").Value = 5004000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8586
$ws.Range("H29").Value = 74004.60000000001
$ws.Range("I29").Value = 69999.5
$ws.Range("K29").Value = 69999.5
$ws.Range("M29").Value = -69709.5
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H113").Value = 1470.8572
$ws.Range("I113").Value = 1181.5
$ws.Range("K113").Value = 1181.5
$ws.Range("M113").Value = 988.5
$ws.Range("H132").Value = 1480713
$ws.Range("J132").Value = 1911.1818
$ws.Range("L132").Value = 5733.5454
$ws.Range("N132").Value = -10793.5454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1821.15
$ws.Range("I132").Value = 1260.3334
$ws.Range("K132").Value = 3781.0002
$ws.Range("M132").Value = -1251.0002
$ws.Range("H136").Value = 3651.5
$ws.Range("I136").Value = 5701
$ws.Range("J136").Value = 2968.3333
$ws.Range("K136").Value = 17103
$ws.Range("L136").Value = 8904.999899999999
$ws.Range("M136").Value = -14553
$ws.Range("N136").Value = -14004.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 659.625
$ws.Range("I107").Value = 488.53845
$ws.Range("K107").Value = 1465.61535
$ws.Range("M107").Value = 454.38465
$ws.Range("H132").Value = 1381.8276
$ws.Range("I132").Value = 978.1951
$ws.Range("K132").Value = 2934.5853
$ws.Range("M132").Value = -404.5853000000002
